$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 19, shifting existing rows 19-21 down to 20-22
$ws.Rows.Item(19).Insert()

# Copy style of date cell (s="2") from the row below (the original row 19, now row 20) into new D19
$ws.Cells.Item(20, 4).Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 19 with the data for the new record
$ws.Cells.Item(19, 1).Value = 7
$ws.Cells.Item(19, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(19, 3).Value = "Ñuble"
$ws.Cells.Item(19, 4).Value = 44784
$ws.Cells.Item(19, 5).Value = 16
$ws.Cells.Item(19, 6).Value = 100112037
$ws.Cells.Item(19, 7).Value = "Cebollín"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 160
$ws.Cells.Item(19, 11).Value = 8000
$ws.Cells.Item(19, 12).Value = 8500
$ws.Cells.Item(19, 13).Value = 8250
$ws.Cells.Item(19, 14).Value = "$/docena de atados"
$ws.Cells.Item(19, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(19, 16).Value = 2750
$ws.Cells.Item(19, 17).Value = 3
$ws.Cells.Item(19, 18).Value = "Hortaliza"
